# Updated cryptos list - apply new Price / Volume(1h) values,
# and correct two swapped Coin/Link pairs (rows 16/17 and 47/48)
# and one Coin/Link replacement (row 51: THORChain -> Monero).
#
# Note: several "Price" (column D) values look like plain numbers
# (e.g. "561.48", "0.641", "1.00"), but in this sheet that column is
# stored as TEXT (to preserve trailing zeros / exact formatting, and
# because other rows use a dotted-thousands style like "66.558.51"
# that can't be a number anyway). Setting NumberFormat to "@" (Text)
# before assigning the value keeps Excel from silently reinterpreting
# these strings as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.558.51"
$ws.Range("E2").Value = "  +3.17%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.515.90"
$ws.Range("E3").Value = "  +6.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.45%  "

# Row 5 - BNB
Set-TextValue "D5" "561.48"
$ws.Range("E5").Value = "  +6.50%  "

# Row 6 - Solana
Set-TextValue "D6" "184.97"
$ws.Range("E6").Value = "  +6.96%  "

# Row 7 - XRP
Set-TextValue "D7" "0.641"
$ws.Range("E7").Value = "  +9.75%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.510.56"
$ws.Range("E8").Value = "  +4.95%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.10%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.641"
$ws.Range("E10").Value = "  +6.41%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.157"
$ws.Range("E11").Value = "  +18.11%  "

# Row 12 - Avalanche
Set-TextValue "D12" "55.40"
$ws.Range("E12").Value = "  +5.16%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000281"
$ws.Range("E13").Value = "  +9.09%  "

# Row 14 - Polkadot
Set-TextValue "D14" "9.41"
$ws.Range("E14").Value = "  +5.77%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.086.30"
$ws.Range("E15").Value = "  +7.77%  "

# Row 16 & 17 - swap Chainlink / WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.514.84"
$ws.Range("E16").Value = "  +7.33%  "

$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D17" "18.70"
$ws.Range("E17").Value = "  +7.62%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  +4.45%  "

# Row 19 - WrappedBTC
$ws.Range("D19").Value = "66.601.98"
$ws.Range("E19").Value = "  +4.21%  "

# Row 20 - Uniswap
Set-TextValue "D20" "12.12"
$ws.Range("E20").Value = "  +8.71%  "

# Row 21 - Polygon
Set-TextValue "D21" "1.00"
$ws.Range("E21").Value = "  +5.13%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "420.61"
$ws.Range("E22").Value = "  +11.32%  "

# Row 23 - PancakeSwap
Set-TextValue "D23" "4.10"
$ws.Range("E23").Value = "  +11.31%  "

# Row 24 - Litecoin
Set-TextValue "D24" "86.69"
$ws.Range("E24").Value = "  +6.72%  "

# Row 25 - Toncoin
Set-TextValue "D25" "4.15"
$ws.Range("E25").Value = "  +0.31%  "

# Row 26 - ImmutableX
Set-TextValue "D26" "2.95"
$ws.Range("E26").Value = "  +9.46%  "

# Row 27 - RenderToken
Set-TextValue "D27" "11.00"
$ws.Range("E27").Value = "  -1.33%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "12.42"
$ws.Range("E28").Value = "  +10.98%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -1.07%  "

# Row 30 - Filecoin
Set-TextValue "D30" "9.15"
$ws.Range("E30").Value = "  +13.31%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "30.46"
$ws.Range("E31").Value = "  +6.62%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "6.85"
$ws.Range("E32").Value = "  +5.02%  "

# Row 33 - Bittensor
Set-TextValue "D33" "628.26"
$ws.Range("E33").Value = "  +0.33%  "

# Row 34 - Cosmos
Set-TextValue "D34" "11.90"
$ws.Range("E34").Value = "  +6.74%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.112"
$ws.Range("E35").Value = "  +7.37%  "

# Row 36 - OKB
Set-TextValue "D36" "60.47"
$ws.Range("E36").Value = "  +7.06%  "

# Row 37 - Kaspa
Set-TextValue "D37" "0.149"
$ws.Range("E37").Value = "  +19.48%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0816"
$ws.Range("E38").Value = "  +10.85%  "

# Row 39 - InjectiveProtocol
Set-TextValue "D39" "38.29"
$ws.Range("E39").Value = "  +6.67%  "

# Row 40 - Dai
Set-TextValue "D40" "0.999"
$ws.Range("E40").Value = "  -0.08%  "

# Row 41 - TheGraph
Set-TextValue "D41" "0.388"
$ws.Range("E41").Value = "  +2.97%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +6.76%  "

# Row 43 - Maker
$ws.Range("D43").Value = "3.105.22"
$ws.Range("E43").Value = "  +8.51%  "

# Row 44 - FirstDigitalUSD
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  +1.21%  "

# Row 45 - Fetch.AI
Set-TextValue "D45" "2.63"
$ws.Range("E45").Value = "  +1.15%  "

# Row 46 - ThetaToken
Set-TextValue "D46" "2.87"
$ws.Range("E46").Value = "  +10.57%  "

# Row 47 & 48 - swap ApeXProtocol / VeChain
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D47" "0.0420"
$ws.Range("E47").Value = "  +6.12%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D48" "3.30"
$ws.Range("E48").Value = "  +9.21%  "

# Row 49 - WEMIXToken
$ws.Range("E49").Value = "  +1.98%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +7.69%  "

# Row 51 - THORChain -> Monero
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D51" "140.17"
$ws.Range("E51").Value = "  +1.95%  "
